# Automatic update of files.
#
# The underlying edit rotates record data among rows 2, 4, 8, 7, 6 and 10
# of the active sheet, in that cyclic order (each row ends up holding the
# data that used to belong to the next row in the cycle):
#   row 2  <- old row 4
#   row 4  <- old row 8
#   row 8  <- old row 7
#   row 7  <- old row 6
#   row 6  <- old row 10
#   row 10 <- old row 2
#
# Only the columns whose value actually differs between two rows in the
# cycle need to be written (A, B, D, E, F, G, H, I, Q, R, AO); every other
# column already holds identical data in both rows. Restricting the
# writes this way also avoids Excel's automatic re-typing of untouched
# text cells (e.g. the date-like strings in columns Y/AA) that would
# otherwise occur if whole rows were blindly copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) that can differ between rows and therefore need to
# be captured before the rotation starts overwriting cells.
$cols = @("A","B","D","E","F","G","H","I","Q","R","AO")

function Snapshot-Row($ws, $row, $cols) {
    $snap = @{}
    foreach ($col in $cols) {
        $snap[$col] = $ws.Range("$col$row").Value2
    }
    return $snap
}

# Snapshot every row involved in the rotation *before* any writes happen.
$row2  = Snapshot-Row $ws 2  $cols
$row4  = Snapshot-Row $ws 4  $cols
$row6  = Snapshot-Row $ws 6  $cols
$row7  = Snapshot-Row $ws 7  $cols
$row8  = Snapshot-Row $ws 8  $cols
$row10 = Snapshot-Row $ws 10 $cols

function Write-Row($ws, $row, $cols, $snap) {
    foreach ($col in $cols) {
        $val = $snap[$col]
        $target = $ws.Range("$col$row")
        if ($val -eq $null -or $val -eq "") {
            $target.Value = ""
        } elseif ($col -eq "I") {
            # Column I ("Antal") is stored as text in the source data
            # (e.g. "30", "7"); force text formatting so Excel's
            # automatic type inference doesn't turn it into a number.
            $target.NumberFormat = "@"
            $target.Value = [string]$val
        } else {
            $target.Value = $val
        }
    }
}

Write-Row $ws 2  $cols $row4
Write-Row $ws 4  $cols $row8
Write-Row $ws 8  $cols $row7
Write-Row $ws 7  $cols $row6
Write-Row $ws 6  $cols $row10
Write-Row $ws 10 $cols $row2
